$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at the top of the data (row 2), pushing all existing
# case rows down by one (row 2 -> 3, 3 -> 4, ... 10 -> 11).
$ws.Rows.Item(2).Insert()

# New timestamp for this scrape pass, applied to every data row (2..11).
$timestamp = "2025-09-27 18:23:03"
for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 1).Value = $timestamp
}

# Fill in the brand-new case that landed in row 2.
$ws.Range("B2").Value = "競馬AIの開発ができる方、もしくはすでに開発済みの方"
$ws.Range("C2").Value = "システム開発"
$ws.Range("D2").Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("E2").Value = "期限情報なし"
$ws.Range("F2").Value = "https://www.lancers.jp/work/detail/5401880"
$ws.Range("F2").Style = $ws.Range("F3").Style
$ws.Range("G2").Value = 375
$ws.Range("H2").Value = "🔥AI,Ai ◆開発"

# The row-insert does not cleanly renumber the hyperlink collection (it
# keeps stale ref/id pairings), so rebuild it from scratch to match the
# shifted F-column URLs plus the new row's link.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5401880")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5401688")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5401604")
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5401806")
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.lancers.jp/work/detail/5401602")
$ws.Hyperlinks.Add($ws.Range("F7"), "https://www.lancers.jp/work/detail/5401800")
$ws.Hyperlinks.Add($ws.Range("F8"), "https://www.lancers.jp/work/detail/5399347")
$ws.Hyperlinks.Add($ws.Range("F9"), "https://www.lancers.jp/work/detail/5401736")
$ws.Hyperlinks.Add($ws.Range("F10"), "https://www.lancers.jp/work/detail/5401572")
$ws.Hyperlinks.Add($ws.Range("F11"), "https://www.lancers.jp/work/detail/5401534")
